# Mark additional "Done [yes or no]" cells as solved ("yes") in column C,
# matching the author's commit "solved some more questions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$doneCells = @(
    "C16", "C22", "C23", "C24", "C26", "C32", "C33", "C35",
    "C62", "C63", "C70", "C97", "C139", "C261", "C296",
    "C410", "C444", "C452"
)

foreach ($addr in $doneCells) {
    $ws.Range($addr).Value = "yes"
}

# Reflect the author's scroll position / active cell at the time of saving.
$ws.Range("C452").Select()
$excel.ActiveWindow.ScrollRow = 439
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Updated $($doneCells.Count) cells to 'yes'"
